$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 52.40544891357422
$ws.Range("C2").Value = -8.547795295715332
$ws.Range("B3").Value = 46.81840896606445
$ws.Range("C3").Value = -6.757501602172852
$ws.Range("B4").Value = 34.35824203491211
$ws.Range("C4").Value = -3.568210124969482
$ws.Range("B5").Value = 17.29892349243164
$ws.Range("C5").Value = -0.4281232357025146
$ws.Range("B6").Value = 14.36217784881592
$ws.Range("C6").Value = 0.01025795936584473
$ws.Range("B7").Value = 13.10221481323242
$ws.Range("C7").Value = 0.1583518385887146
$ws.Range("B8").Value = 11.9117431640625
$ws.Range("C8").Value = 0.2958714365959167
$ws.Range("B9").Value = 11.88819599151611
$ws.Range("C9").Value = 0.2927542924880981
$ws.Range("B10").Value = 11.46057510375977
$ws.Range("C10").Value = 0.3718902468681335
$ws.Range("B11").Value = 10.67667865753174
$ws.Range("C11").Value = 0.420953094959259
$ws.Range("B12").Value = 10.54313182830811
$ws.Range("C12").Value = 0.4576348662376404
$ws.Range("B13").Value = 10.17471218109131
$ws.Range("C13").Value = 0.4804003238677979
$ws.Range("B14").Value = 9.824414253234863
$ws.Range("C14").Value = 0.5154540538787842
$ws.Range("B15").Value = 9.503444671630859
$ws.Range("C15").Value = 0.5291653871536255
$ws.Range("B16").Value = 9.515039443969727
$ws.Range("C16").Value = 0.5458121299743652
$ws.Range("B17").Value = 8.62856388092041
$ws.Range("C17").Value = 0.6253806352615356
$ws.Range("B18").Value = 8.494733810424805
$ws.Range("C18").Value = 0.641605019569397
$ws.Range("B19").Value = 8.588038444519043
$ws.Range("C19").Value = 0.6405611634254456
$ws.Range("B20").Value = 8.093818664550781
$ws.Range("C20").Value = 0.6663030385971069
$ws.Range("B21").Value = 8.45398998260498
$ws.Range("C21").Value = 0.6373758316040039
$ws.Range("B22").Value = 8.554228782653809
$ws.Range("C22").Value = 0.6307522058486938
$ws.Range("B23").Value = 7.638081073760986
$ws.Range("C23").Value = 0.6969456672668457
$ws.Range("B24").Value = 7.551326274871826
$ws.Range("C24").Value = 0.7162183523178101
$ws.Range("B25").Value = 7.465870380401611
$ws.Range("C25").Value = 0.7159451842308044
$ws.Range("B26").Value = 7.51484489440918
$ws.Range("C26").Value = 0.7141540050506592
$ws.Range("B27").Value = 7.385990142822266
$ws.Range("C27").Value = 0.7066656351089478
$ws.Range("B28").Value = 7.148406982421875
$ws.Range("C28").Value = 0.7444168329238892
$ws.Range("B29").Value = 7.007951259613037
$ws.Range("C29").Value = 0.7345058917999268
$ws.Range("B30").Value = 7.092115879058838
$ws.Range("C30").Value = 0.7381138205528259
$ws.Range("B31").Value = 7.044261932373047
$ws.Range("C31").Value = 0.7389810085296631
$ws.Range("B32").Value = 7.078729629516602
$ws.Range("C32").Value = 0.7300344705581665
$ws.Range("B33").Value = 6.744141578674316
$ws.Range("C33").Value = 0.7561941146850586
$ws.Range("B34").Value = 6.728821754455566
$ws.Range("C34").Value = 0.7592419981956482
$ws.Range("B35").Value = 6.586283206939697
$ws.Range("C35").Value = 0.780877947807312
$ws.Range("B36").Value = 6.427896499633789
$ws.Range("C36").Value = 0.7780163884162903
$ws.Range("B37").Value = 6.432909965515137
$ws.Range("C37").Value = 0.7741916179656982
$ws.Range("B38").Value = 6.603747844696045
$ws.Range("C38").Value = 0.7715558409690857
$ws.Range("B39").Value = 6.33154821395874
$ws.Range("C39").Value = 0.7899177670478821
$ws.Range("B40").Value = 6.084758758544922
$ws.Range("C40").Value = 0.8082225322723389
$ws.Range("B41").Value = 6.289232730865479
$ws.Range("C41").Value = 0.7875809669494629
$ws.Range("B42").Value = 6.236540794372559
$ws.Range("C42").Value = 0.7956802248954773
$ws.Range("B43").Value = 6.089844703674316
$ws.Range("C43").Value = 0.8005338907241821
$ws.Range("B44").Value = 6.340396404266357
$ws.Range("C44").Value = 0.7857683897018433
$ws.Range("B45").Value = 6.036722660064697
$ws.Range("C45").Value = 0.8039543628692627
$ws.Range("B46").Value = 5.884133338928223
$ws.Range("C46").Value = 0.8104521632194519
$ws.Range("B47").Value = 5.972614288330078
$ws.Range("C47").Value = 0.8177534341812134
$ws.Range("B48").Value = 5.968766212463379
$ws.Range("C48").Value = 0.8072613477706909
$ws.Range("B49").Value = 5.956734657287598
$ws.Range("C49").Value = 0.8082748651504517
$ws.Range("B50").Value = 6.048097133636475
$ws.Range("C50").Value = 0.8112385869026184
$ws.Range("B51").Value = 5.791770458221436
$ws.Range("C51").Value = 0.8219859600067139
$ws.Range("B52").Value = 5.844794273376465
$ws.Range("C52").Value = 0.8207231760025024
$ws.Range("B53").Value = 5.726579666137695
$ws.Range("C53").Value = 0.8211728930473328
$ws.Range("B54").Value = 5.759580135345459
$ws.Range("C54").Value = 0.8164608478546143
$ws.Range("B55").Value = 5.86553430557251
$ws.Range("C55").Value = 0.8184210062026978
$ws.Range("B56").Value = 5.660834789276123
$ws.Range("C56").Value = 0.8300696611404419
$ws.Range("B57").Value = 5.815227031707764
$ws.Range("C57").Value = 0.819373607635498
$ws.Range("B58").Value = 5.543291091918945
$ws.Range("C58").Value = 0.8376485705375671
$ws.Range("B59").Value = 5.49190616607666
$ws.Range("C59").Value = 0.8368971347808838
$ws.Range("B60").Value = 5.500041484832764
$ws.Range("C60").Value = 0.8446533679962158
$ws.Range("B61").Value = 5.811877727508545
$ws.Range("C61").Value = 0.8048028349876404
$ws.Range("B62").Value = 5.399606704711914
$ws.Range("C62").Value = 0.8432849645614624
$ws.Range("B63").Value = 5.826673984527588
$ws.Range("C63").Value = 0.8180335164070129
$ws.Range("B64").Value = 5.551222324371338
$ws.Range("C64").Value = 0.8374972343444824
$ws.Range("B65").Value = 5.506695747375488
$ws.Range("C65").Value = 0.8374452590942383
$ws.Range("B66").Value = 5.489109039306641
$ws.Range("C66").Value = 0.8298180103302002
$ws.Range("B67").Value = 5.504336357116699
$ws.Range("C67").Value = 0.8374220728874207
$ws.Range("B68").Value = 5.505527019500732
$ws.Range("C68").Value = 0.8307868838310242
$ws.Range("B69").Value = 5.322421073913574
$ws.Range("C69").Value = 0.8502706289291382
$ws.Range("B70").Value = 5.484179496765137
$ws.Range("C70").Value = 0.8451113700866699
$ws.Range("B71").Value = 5.645872116088867
$ws.Range("C71").Value = 0.8273598551750183
$ws.Range("B72").Value = 5.427425384521484
$ws.Range("C72").Value = 0.8432795405387878
$ws.Range("B73").Value = 5.218438148498535
$ws.Range("C73").Value = 0.8524907231330872
$ws.Range("B74").Value = 5.426321983337402
$ws.Range("C74").Value = 0.8465814590454102
$ws.Range("B75").Value = 5.373984813690186
$ws.Range("C75").Value = 0.8444929122924805
$ws.Range("B76").Value = 5.129755973815918
$ws.Range("C76").Value = 0.859188437461853
$ws.Range("B77").Value = 5.152495861053467
$ws.Range("C77").Value = 0.8567067384719849
$ws.Range("B78").Value = 5.28340482711792
$ws.Range("C78").Value = 0.8545277714729309
$ws.Range("B79").Value = 5.439005851745605
$ws.Range("C79").Value = 0.8456141948699951
$ws.Range("B80").Value = 5.321115016937256
$ws.Range("C80").Value = 0.8465718626976013
$ws.Range("B81").Value = 5.30781078338623
$ws.Range("C81").Value = 0.8530145287513733
$ws.Range("B82").Value = 5.275221824645996
$ws.Range("C82").Value = 0.8498678207397461
$ws.Range("B83").Value = 5.360738754272461
$ws.Range("C83").Value = 0.8512598276138306
$ws.Range("B84").Value = 5.254478931427002
$ws.Range("C84").Value = 0.842096745967865
$ws.Range("B85").Value = 5.340359210968018
$ws.Range("C85").Value = 0.8410004377365112
$ws.Range("B86").Value = 5.078260898590088
$ws.Range("C86").Value = 0.8642523884773254
$ws.Range("B87").Value = 5.075554370880127
$ws.Range("C87").Value = 0.8659154176712036
$ws.Range("B88").Value = 5.217664241790771
$ws.Range("C88").Value = 0.8507083654403687
$ws.Range("B89").Value = 5.018243789672852
$ws.Range("C89").Value = 0.86647629737854
$ws.Range("B90").Value = 5.118835926055908
$ws.Range("C90").Value = 0.8565160632133484
$ws.Range("B91").Value = 5.136265277862549
$ws.Range("C91").Value = 0.8609759211540222
$ws.Range("B92").Value = 5.101777076721191
$ws.Range("C92").Value = 0.8664219379425049
$ws.Range("B93").Value = 5.183019161224365
$ws.Range("C93").Value = 0.8622119426727295
$ws.Range("B94").Value = 4.970968723297119
$ws.Range("C94").Value = 0.8668940663337708
$ws.Range("B95").Value = 5.068518161773682
$ws.Range("C95").Value = 0.8643868565559387
$ws.Range("B96").Value = 4.923008918762207
$ws.Range("C96").Value = 0.8696881532669067
$ws.Range("B97").Value = 5.269438743591309
$ws.Range("C97").Value = 0.8547347784042358
$ws.Range("B98").Value = 5.007972717285156
$ws.Range("C98").Value = 0.8675112128257751
$ws.Range("B99").Value = 5.032581806182861
$ws.Range("C99").Value = 0.8707244396209717
$ws.Range("B100").Value = 4.917586326599121
$ws.Range("C100").Value = 0.8734118938446045
$ws.Range("B101").Value = 5.083066940307617
$ws.Range("C101").Value = 0.8963643801361084
